# "Generate Report for Handback" - update the localization-status workbook
# to reflect that both files have now been handed back and are in sync
# with en-US.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

$a46Url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/41efd7644204d35d32fa1454c634e3a3e4c32a3b/e2e/a46d0e08-2453-416e-9b29-e3a1a91e3d7a.md"
$aacUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/41efd7644204d35d32fa1454c634e3a3e4c32a3b/e2e/aac79222-062a-4681-b8ed-d0285ad8e595.md"

$a46Name = "a46d0e08-2453-416e-9b29-e3a1a91e3d7a.md"
$aacName = "aac79222-062a-4681-b8ed-d0285ad8e595.md"

# --- Overview sheet: status cells move from "Ready for handoff" to handed back ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText

# widen the status columns to fit the longer text
$wsOverview.Columns.Item(5).ColumnWidth = 29.083333333333332
$wsOverview.Columns.Item(6).ColumnWidth = 29.083333333333332

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $statusText
$wsZh.Range("C3").Value = $statusText

# row 2 (a46d0e08...)
$wsZh.Range("I2").Value = $a46Name
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $a46Url, "", "", $a46Name)
$wsZh.Range("J2").Value = "a46d0e08-2453-416e-9b29-e3a1a91e3d7a.263c7478941e6f4247ab41b71d5edd3aa115a34b.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-16 06:55:01"

# row 3 (aac79222...)
$wsZh.Range("I3").Value = $aacName
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $aacUrl, "", "", $aacName)
$wsZh.Range("J3").Value = "aac79222-062a-4681-b8ed-d0285ad8e595.e42f4eab8cc73b3badd248a26e0ded876ad41f8d.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-08-16 06:55:01"

$wsZh.Columns.Item(3).ColumnWidth = 29.083333333333332
$wsZh.Columns.Item(9).ColumnWidth = 39.083333333333336
$wsZh.Columns.Item(10).ColumnWidth = 39.083333333333336

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $statusText
$wsDe.Range("C3").Value = $statusText

# row 2 (a46d0e08...)
$wsDe.Range("I2").Value = $a46Name
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $a46Url, "", "", $a46Name)
$wsDe.Range("J2").Value = "a46d0e08-2453-416e-9b29-e3a1a91e3d7a.263c7478941e6f4247ab41b71d5edd3aa115a34b.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-16 06:55:15"

# row 3 (aac79222...)
$wsDe.Range("I3").Value = $aacName
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $aacUrl, "", "", $aacName)
$wsDe.Range("J3").Value = "aac79222-062a-4681-b8ed-d0285ad8e595.e42f4eab8cc73b3badd248a26e0ded876ad41f8d.de-de.xlf"
$wsDe.Range("K3").Value = "2016-08-16 06:55:15"

$wsDe.Columns.Item(3).ColumnWidth = 29.083333333333332
$wsDe.Columns.Item(9).ColumnWidth = 39.083333333333336
$wsDe.Columns.Item(10).ColumnWidth = 39.083333333333336
